$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Total_Citations_(2023" row (row 23) text values so the
# range minimums reflect the newly-zeroed "Unknown" counts.
$ws.Range("B23").Value = "363 (range: 0 to 671)"
$ws.Range("C23").Value = "1,528 (range: 0 to 10,260)"
$ws.Range("D23").Value = "192 (range: 0 to 2,672)"
$ws.Range("E23").Value = "215 (range: 0 to 644)"
$ws.Range("F23").Value = "334 (range: 0 to 3,163)"
$ws.Range("G23").Value = "10,773 (range: 0 to 169,257)"
$ws.Range("H23").Value = "4,133 (range: 0 to 15,604)"
$ws.Range("I23").Value = "8,339 (range: 0 to 102,352)"

# Remove the "Unknown" count row (row 24) entirely - its counts were
# zeroed out so the row is dropped and everything below shifts up.
$ws.Rows.Item(24).Delete()
